$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pipeline incidence matrix values that changed between the
# original layout and the revised (robust) pipeline design.

$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("I2").Value = 1

$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0

$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("I4").Value = 1

$ws.Range("C5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("H5").Value = 1

$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 1

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("F7").Value = 1

$ws.Range("C8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("I8").Value = 0

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("H9").Value = 0

# Mirror the cursor/selection ending on I2, matching the saved view state.
$ws.Range("I2").Select()
